$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = @{ B = 12766.6;            C = 13522.66666666667 }
    4  = @{ B = 12333.2;            C = 13845.33333333333 }
    5  = @{ B = 11899.8;            C = 14168 }
    6  = @{ B = 11466.4;            C = 14490.66666666667 }
    7  = @{ B = 14217.34285714286;  C = 15056.17142857143 }
    8  = @{ B = 16968.28571428571;  C = 15621.67619047619 }
    9  = @{ B = 19719.22857142857;  C = 16187.18095238095 }
    10 = @{ B = 22470.17142857143;  C = 16752.68571428572 }
    11 = @{ B = 11299.2;            C = 9820.800000000001 }
    12 = @{ B = 11404.8;            C = 7444.799999999999 }
    13 = @{ B = 13886.4;            C = 13992 }
    14 = @{ B = 12196.8;            C = 6969.6 }
    15 = @{ B = 10771.2;            C = 8764.799999999999 }
    16 = @{ B = 7920;               C = 15364.8 }
    17 = @{ B = 15787.2;            C = 5755.200000000001 }
    18 = @{ B = 14361.6;            C = 11668.8 }
    19 = @{ B = 9292.799999999999;  C = 8289.6 }
    20 = @{ B = 22968;              C = 15364.8 }
    21 = @{ B = 23073.6;            C = 17529.6 }
    22 = @{ B = 17265.6;            C = 18691.2 }
    23 = @{ B = 26188.8;            C = 20803.2 }
    24 = @{ B = 22334.4;            C = 11035.2 }
    25 = @{ B = 15892.8;            C = 11985.6 }
    26 = @{ B = 25766.4;            C = 13780.8 }
    27 = @{ B = 19694.4;            C = 20064 }
    28 = @{ B = 8606.4;             C = 21648 }
    29 = @{ B = 9873.6;             C = 16526.4 }
    30 = @{ B = 11457.6;            C = 24763.2 }
    31 = @{ B = 6916.8;             C = 18216 }
    32 = @{ B = 9662.4;             C = 16209.6 }
    33 = @{ B = 10507.2;            C = 24499.2 }
    34 = @{ B = 12249.6;            C = 21489.6 }
    35 = @{ B = 14308.8;            C = 17424 }
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row].B
    $ws.Range("C$row").Value = $values[$row].C
}
